# Generate Report for Handback
# Update the "Correspond Handoff Datetime" (col E) and
# "Correspond Handback DateTime" (col H) values for the fc6afe69... file row
# (row 3) on both the zh-cn and de-de sheets, reflecting a newer handback run.

$wb = $excel.ActiveWorkbook

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("E3").Value = "2016-03-22 22:55:20"
$zhcn.Range("H3").Value = "2016-03-22 22:55:46"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("E3").Value = "2016-03-22 22:55:26"
$dede.Range("H3").Value = "2016-03-22 22:55:52"
